# Refactor psf_beads_schema.yaml: rename PSFBeadsKeyValues -> PSFBeadsKeyMeasurements,
# rename the key_values column reference on PSFBeadsOutput to key_measurements, and
# rebuild the PSFBeadsKeyMeasurements header row as a (wider) table of measurements.

$wb = $excel.ActiveWorkbook

# 1) Rename the sheet tab itself.
$wsKeyValues = $wb.Worksheets.Item("PSFBeadsKeyValues")
$wsKeyValues.Name = "PSFBeadsKeyMeasurements"

# 2) On PSFBeadsOutput, the column that used to point at "key_values" now points at
#    "key_measurements" (column I).
$wsOutput = $wb.Worksheets.Item("PSFBeadsOutput")
$wsOutput.Range("I1").Value = "key_measurements"

# 3) Rebuild the header row of PSFBeadsKeyMeasurements with the new/renamed/reordered
#    columns. The table grew from 43 columns (A:AQ) to 54 columns (A:BB).
$wsKeyMeasurements = $wb.Worksheets.Item("PSFBeadsKeyMeasurements")

$headers = @(
    "channel_nr",
    "considered_valid_count",
    "considered_self_proximity_count",
    "considered_lateral_edge_count",
    "considered_axial_edge_count",
    "considered_intensity_outlier_count",
    "considered_bad_z_fit_count",
    "considered_bad_y_fit_count",
    "considered_bad_x_fit_count",
    "intensity_max_mean",
    "intensity_max_median",
    "intensity_max_std",
    "intensity_min_mean",
    "intensity_min_median",
    "intensity_min_std",
    "intensity_std_mean",
    "intensity_std_median",
    "intensity_std_std",
    "fit_r2_z_mean",
    "fit_r2_z_median",
    "fit_r2_z_std",
    "fit_r2_y_mean",
    "fit_r2_y_median",
    "fit_r2_y_std",
    "fit_r2_x_mean",
    "fit_r2_x_median",
    "fit_r2_x_std",
    "fwhm_pixel_z_mean",
    "fwhm_pixel_z_median",
    "fwhm_pixel_z_std",
    "fwhm_pixel_y_mean",
    "fwhm_pixel_y_median",
    "fwhm_pixel_y_std",
    "fwhm_pixel_x_mean",
    "fwhm_pixel_x_median",
    "fwhm_pixel_x_std",
    "fwhm_micron_z_mean",
    "fwhm_micron_z_median",
    "fwhm_micron_z_std",
    "fwhm_micron_y_mean",
    "fwhm_micron_y_median",
    "fwhm_micron_y_std",
    "fwhm_micron_x_mean",
    "fwhm_micron_x_median",
    "fwhm_micron_x_std",
    "fwhm_lateral_asymmetry_ratio_mean",
    "fwhm_lateral_asymmetry_ratio_median",
    "fwhm_lateral_asymmetry_ratio_std",
    "column_series",
    "table_data",
    "data_reference",
    "linked_references",
    "name",
    "description"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsKeyMeasurements.Cells.Item(1, $i + 1).Value = $headers[$i]
}
